$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.4778628081766
$ws.Range("C2").Value = 17.1175986105978
$ws.Range("D2").Value = 17.8381270057554
$ws.Range("C3").Value = 18.1092138950899
$ws.Range("D3").Value = 19.2394096308851
$ws.Range("B11").Value = 13.9286135866664
$ws.Range("C11").Value = 13.2333145854194
$ws.Range("D11").Value = 14.6239125879133
$ws.Range("B13").Value = 17.3020839614069
$ws.Range("C13").Value = 16.6285285100235
$ws.Range("D13").Value = 17.9756394127904
$ws.Range("B14").Value = 25.8176804364833
$ws.Range("C14").Value = 25.4207830126572
$ws.Range("D14").Value = 26.2145778603093
$ws.Range("C15").Value = 28.0556145330742
$ws.Range("D15").Value = 29.2819146425341
$ws.Range("B23").Value = 24.9676333375889
$ws.Range("C23").Value = 24.0996816084345
$ws.Range("D23").Value = 25.8355850667434
$ws.Range("B25").Value = 23.1977463905367
$ws.Range("C25").Value = 22.4841233151793
$ws.Range("D25").Value = 23.9113694658941
$ws.Range("B26").Value = 26.888852791522
$ws.Range("C26").Value = 26.2840599100719
$ws.Range("D26").Value = 27.4936456729721
$ws.Range("C27").Value = 26.9617128666953
$ws.Range("D27").Value = 28.7918085096094
$ws.Range("B35").Value = 23.7925359976382
$ws.Range("C35").Value = 22.4533074894758
$ws.Range("D35").Value = 25.1317645058005
$ws.Range("B37").Value = 26.3869099551691
$ws.Range("C37").Value = 25.2573013744177
$ws.Range("D37").Value = 27.5165185359205
$ws.Range("B38").Value = 18.5840386333021
$ws.Range("C38").Value = 17.9671758541937
$ws.Range("D38").Value = 19.2009014124105
$ws.Range("C39").Value = 18.6186128419952
$ws.Range("D39").Value = 20.4680649318187
$ws.Range("B47").Value = 14.7145599857929
$ws.Range("C47").Value = 13.4191286852481
$ws.Range("D47").Value = 16.0099912863377
$ws.Range("B49").Value = 18.8176682846185
$ws.Range("C49").Value = 17.6445715829381
$ws.Range("D49").Value = 19.9907649862988
$ws.Range("B50").Value = 18.2790863468141
$ws.Range("C50").Value = 17.7086767334779
$ws.Range("D50").Value = 18.8494959601504
$ws.Range("C51").Value = 17.8436710027609
$ws.Range("D51").Value = 19.5229413950389
$ws.Range("B59").Value = 17.7049808040103
$ws.Range("C59").Value = 16.2815627940052
$ws.Range("D59").Value = 19.1283988140154
$ws.Range("B61").Value = 17.7140126023854
$ws.Range("C61").Value = 16.6641188900468
$ws.Range("D61").Value = 18.763906314724
$ws.Range("B62").Value = 16.7223157348188
$ws.Range("C62").Value = 16.1687738455109
$ws.Range("D62").Value = 17.2758576241267
$ws.Range("C63").Value = 15.6153160964197
$ws.Range("D63").Value = 17.2511550460299
$ws.Range("B71").Value = 15.9779234610467
$ws.Range("C71").Value = 14.5525384282777
$ws.Range("D71").Value = 17.4033084938156
$ws.Range("B73").Value = 16.7126947567179
$ws.Range("C73").Value = 15.7084970354776
$ws.Range("D73").Value = 17.7168924779581
$ws.Range("B74").Value = 16.2405917484416
$ws.Range("C74").Value = 15.7116001581493
$ws.Range("D74").Value = 16.7695833387339
$ws.Range("C75").Value = 16.4312908880532
$ws.Range("D75").Value = 17.9770834330384
$ws.Range("B83").Value = 14.0638042866307
$ws.Range("C83").Value = 12.9045274456188
$ws.Range("D83").Value = 15.2230811276426
$ws.Range("B85").Value = 15.5556817545784
$ws.Range("C85").Value = 14.5412204070847
$ws.Range("D85").Value = 16.570143102072
